# Actualización automática 2025-08-06 14:45:08
#
# Updates sales figures for ALMEIDA CUATIN JHONATHANN CARLOS / COMFALASDI
# COMPAÑIA FAMILIAR LASCANO DIAZ C. LTDA. across the three report sheets,
# and refreshes the dependent aggregate rows (counts, totals, deltas and
# completion percentages) that are kept as static computed values in this
# workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" - per product-group sales for the client
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("D8").Value = 457.92
$wsGrupo.Range("M8").Value = 1998.01
$wsGrupo.Range("O8").Value = 1.73
$wsGrupo.Range("R8").Value = 130.5

# Row 33 keeps a textual "<n> de 31" count of clients with non-zero sales
# per product group; the four groups touched above gain one more client.
$wsGrupo.Range("D33").Value = "2 de 31"
$wsGrupo.Range("M33").Value = "2 de 31"
$wsGrupo.Range("O33").Value = "1 de 31"
$wsGrupo.Range("R33").Value = "1 de 31"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" - monthly sales for the client (agosto column)
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F8").Value = 2588.16
$wsMensual.Range("F33").Value = 3451.14

# Column widened slightly (Excel's "best fit" reflow) now that the agosto
# figures contain more digits. (ColumnWidth stores 5/6 of a character less
# than the rendered grid width, so subtract that padding to land on 13.)
$wsMensual.Columns.Item(6).ColumnWidth = 13 - (5/6)

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" - budget vs sales completion by group
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumpl.Range("D3").Value = 915.84
$wsCumpl.Range("E3").Value = 2204.2745
$wsCumpl.Range("F3").Value = 0.2935276894485763

$wsCumpl.Range("D13").Value = 130.5
$wsCumpl.Range("E13").Value = -110.5
$wsCumpl.Range("F13").Value = 6.525

$wsCumpl.Range("D16").Value = 2040.05
$wsCumpl.Range("E16").Value = 19833.05
$wsCumpl.Range("F16").Value = 0.09326752952256424

$wsCumpl.Range("D18").Value = 1.73
$wsCumpl.Range("E18").Value = 1598.27
$wsCumpl.Range("F18").Value = 0.00108125

$wsCumpl.Range("D19").Value = 3451.14
$wsCumpl.Range("E19").Value = 28658.14107555787
$wsCumpl.Range("F19").Value = 0.1074810735213585

# VENTA column widened, CUMPLIMIENTO column narrowed slightly as the
# "best fit" reflow runs over the new figures.
$wsCumpl.Columns.Item(4).ColumnWidth = 13 - (5/6)
$wsCumpl.Columns.Item(6).ColumnWidth = 25 - (5/6)
